$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 15-17); new extent is A1:G14
$ws.Rows("15:17").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = '2023 설 스페셜 혜택전'
$ws.Range("B2").Value = '단 3일간, 최대 3만원 선물세트 전용쿠폰'
$ws.Range("C2").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003554'
$ws.Range("D2").Value = '2023 설 스페셜 혜택전'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2022-12-26'
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '2022-12-28'
$ws.Range("G2").Value = '[''이벤트/쿠폰 > 2023 설 스페셜 혜택전'', ''스마일클럽'', ''4종쿠폰'', ''카드사 청구할인'', ''명절 종합 선물세트 4종 쿠폰'', ''      명절 선물 전용 3,000원 쿠폰 (점포명절, 명절)'', ''      명절 선물 전용 7,000원 쿠폰 (점포명절, 명절)'', ''      명절 선물 전용 15,000원 쿠폰 (점포명절, 명절)'', ''      명절 선물 전용 30,000원 쿠폰 (점포명절, 명절)'', ''쿠폰 발급 후 [My SSG] ▶ [쿠폰] ▶ [보유쿠폰] 에서 확인 가능합니다.'', ''[명절] [점포 명절] 엠블렘이 있는 설 명절 기프트 상품 전용 쿠폰입니다.'', ''[카드 50% 할인]LG 선물 세트)리튠 녹용진액 골드 3종 기획세트'', ''[SSGPAY머니 증정][1/5(목)이후 순차출고][냉동][오늘의소] 꽃갈비를 품은 블랙앵거스 LA갈비 선물세트 3kg'', ''[20만원↑구매시3만원쓱머니증정][배송일 해피콜][직접배송][SSG정담][과일愛]프라이빗 과일선물세트 후룻디니뜨 혼합 9호 (샤인머스켓3 / 사과2 / 배2 / 자몽2 - 4.8kg) - 보자기포장'', ''[냉장][장흥명품한우할인직판장] 구이1호 1.3kg 외 프리미엄 한우 선물세트 특가'', ''최대 57% 할인! GC녹십자 설 명절 건강 선물세트'', ''[백화점라인/국내배송] 폴로 랄프로렌 이태리산 울 100% 시그니처 멀티로고 폴로 머플러 (6color)'', ''[해외직구] 포켓몬 몬스터볼 변신 피규어 6종세트 무료배송 빠른 출고'', ''풋웨어 BEST 에어맥스97/데이브레이크 ~44%'', ''새해맞이 댕댕이 한복/배변패드/유모차 할인 구경하시개!'', ''설맞이 설레는 고양이 모래/캣타워/스크래쳐 용품 할인전'', ''- 삼성카드 청구할인 -'']'

# Row 3
$ws.Range("A3").Value = '설 선물세트 다함께 공동구매 1탄 _프리미엄'
$ws.Range("B3").Value = '최대 40%할인 + 목표달성하면 SSG머니 최대 6만원 혜택'
$ws.Range("C3").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003545'
$ws.Range("D3").Value = '[명절]설 공동구매 1탄'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2022-12-26'
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '2023-01-01'
$ws.Range("G3").Value = '[''이벤트/쿠폰 > [명절]설 공동구매 1탄'', ''스마일클럽'', ''명절 선물세트 7종 최대 40% 할인 + 공동구매 혜택 최대 6만원 적립'', ''공구 이벤트 기간 : 12월 26일 ~ 1월 1일'', ''[오메드] 오일 비네거 세트 30% 할인 + 최대 8천원 적립!'', ''[꽁피튀르] 파리지엥잼 선물세트 30% 할인 + 최대 8천원 적립!'', ''[SSG 정담] 프리미엄 선물세트 33% 할인 + 최대 8천원 적립!'', ''프리미엄 망고혼합 선물세트 17% 할인 + 최대 1만원 적립!'', ''[종근당건강] 6년근 홍삼정업 세트 20% 할인 + 최대 1만원 적립!'', ''[고려천홍삼] 프리미엄 녹용홍삼 37% 할인 + 최대 8천원 적립!'', ''호주산 LA갈비 선물세트 28% 할인 + 최대 8천원 적립!'', ''본 이벤트에 대한 세부사항은 당사 사정에 따라 임의로 변경 혹은 종료될 수 있습니다.'']'

# Row 4
$ws.Range("A4").Value = 'Happy Holiday 연말엔 선물하기로 쓱-'
$ws.Range("B4").Value = '선물하기 처음이라면 SSG머니 행운까지!'
$ws.Range("C4").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003313'
$ws.Range("D4").Value = '선물하기로 쓱- 마음을 전하세요'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2022-12-15'
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '2023-01-22'
$ws.Range("G4").Value = '[''이벤트/쿠폰 > 선물하기로 쓱- 마음을 전하세요'', ''스마일클럽'', ''첫 선물하기 기념 SSG 머니 페이백 이벤트 바로보기'', ''모바일쿠폰'', ''천명에게 쏟아지는 SSG머니 이벤트'', ''선물하기가 처음이시라면, 쓱머니 이벤트 응모하기'', ''※ 본 이벤트에 대한 세부 사항은 당사 사정에 따라 임의로 변경 혹은 조기 종료될 수 있습니다.'', ''※ 모바일 상품 교환권, 모바일 쿠폰, 외식 식사권 등 e쿠폰 선물하기는 해당 이벤트에서 제외 됩니다.'']'

# Row 5
$ws.Range("A5").Value = '최대 8만원 혜택'
$ws.Range("B5").Value = '+ 쓸 때마다 최대 12% 적립'
$ws.Range("C5").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000000858&siteNo=6005&recruitmentPath=L6007001&eventCode=HPG02'
$ws.Range("D5").Value = 'SSG.COM카드 Edition 2 이벤트 안내 페이지'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2022-07-08'
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '2023-06-28'
$ws.Range("G5").Value = '[''이벤트/쿠폰 > SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''스마일클럽 혜택보기'', ''이달의 혜택 01. 한달 내내 8% 즉시 할인 (대상 상품당, 최대 3만원) 바로보기'', ''이달의 혜택 02. SSG.COM카드 Edition2 첫결제 시 1만원 쿠폰 할인 바로보기'', ''혜택 01. SSG머니 최대 12% 적립 바로보기'', ''혜택 02. 매월 스마일클럽 가입비 3,900원 지원 바로보기'', ''이달의 혜택 01. 한달 내내 8% 즉시할인'', ''SSG.COM에서 SSG.COM카드 Edition2로 상품상세 내 [카드 혜택가 > SSG.COM카드 Edition2 더보기 > 즉시할인가] 대상상품 결제 시 (대상 상품당, 최대 3만원 할인)'', ''8% 즉시할인'', ''② [카드 혜택가 > SSG.COM카드 Edition2 더보기 > 즉시할인가] 대상상품결제 시 적용 (대상 상품당, 최대 3만원 할인)'', ''      본 혜택은 SSG.COM카드 Edition2 즉시할인 상품 결제시에만 적용 가능'', ''      즉시할인 상품 확인 방법 : [상품상세 내 > 카드혜택가 > SSG.COM카드 Edition2 행사 더보기]에서 확인'', ''      이벤트 기간 및 혜택 제공 시점에 SSG.COM카드 Edition2를 정상 보유한 회원에 한해 혜택 제공'', ''이달의 혜택 02. SSG.COM에서 SSG.COM카드 Edition2 첫 결제 시 1만원 할인'', ''      직전 6개월간 (2022년 6월 1일 부터 2022년 11월 31일) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한함'', ''첫 결제 쿠폰 이용 방법'', ''② SSGPAY에 등록된 SSG.COM카드 Edition2로'', ''③ SSG.COM에서 기간 내 1만 1천원 이상 첫 결제 시 1만원 쿠폰 할인'', ''      본 혜택은 SSG.COM카드 및 SSG.COM카드 Edition2 로 SSG.COM(이마트몰, 신세계몰, 신세계백화점몰 등)에서 직전 6개월간 (2022년 6월 1일 부터 2022년 11월 30일 까지) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한해 제공됨'', "      본 혜택은 결제 시 [결제방법 > SSGPAY카드] 내 SSG.COM카드 Edition2 선택 시, ''카드할인 최적가'' 추천에 의해 할인 금액이 자동 적용됨. 단, 1만 1천원 이상 첫 결제 시 적용)", ''다운받은 쿠폰은 SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''쿠폰 사용하러 가기'', ''SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''* 쿠폰할인, 적립금 사용 등 선할인을 제외한 카드 최종 결제금액이 2만원 이상인 경우에 한해 혜택 적용'', ''혜택 01. 장 볼 때마다 SSG머니 최대 12% 적립'', ''SSG.COM에서 최대 12%'', ''       스마일클럽으로 5% 적립'', ''       쓱·새벽·트레이더스 구매 시 (구매 전 스마일클럽 적립 아이콘을 꼭 확인해주세요)'', ''       해당 5% 적립은 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능'', ''       SSG.COM 카드 Edition2로 7% 적립'', ''어디서나 한도없이 0.5%'', ''SSG.COM카드 Edition2로 어디서나 한도없이 0.5% 적립 (SSG.COM 외 모든 가맹점)'', ''SSG Money 최대 12% 적립'', ''스마일클럽 5% 적립 + SSG.COM카드 Edition2 최대 7% 적립'', ''      스마일클럽 5% 적립은 쓱·새벽배송·트레이더스 이용 시에 한함'', ''      SSG.COM카드 Edition2 최대 7% 적립'', ''      SSG.COM에서 결제 시 7%(1만 쓱머니 한도), 그 외 가맹점 0.5% 적립(적립한도 제한 없음)'', ''      무이자 할부 및 현대카드에서 제공하는 다른 할인 서비스 이용 시 적립 제외'', ''스마일클럽 가입비 3,900원 매월 100% 지원'', ''      1. 스마일클럽 자동 가입에'', ''      2. SSGPAY 내 카드 자동 등록'', ''      3. 스마일클럽 정기결제수단 자동 등록 및 월 이용료 3,900원 지원까지!'', ''SSG.COM카드 Edition2를 스마일클럽 월 정기결제 수단에 등록 및 전월 이용금액 30만원 이상 시 혜택 제공'', ''TIP. 스마일클럽 가입 시 SSG.COM 혜택'', ''(SSG 가입 시) 15,000원만 담아도 쓱 · 새벽배송 무료배송'', ''쓱 · 새벽배송 · 트레이더스 SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''스마일클럽 단독 혜택 스타벅스 상품 전용 딜'', ''스마일클럽 가입비 매월 3,900원 지원'', ''월 1회, 매달 스마일클럽 정기결제일에 혜택 제공'', ''SSG.COM카드 Edition2는 최초 발급 시, 스마일클럽 월 정기결제 수단에 자동 등록 됨'', ''전월 이용금액 30만원 미만 시, SSG.COM카드 Edition2로 스마일클럽 정기 결제 금액이 자동 결제됨'', ''스마일클럽 무료 이용 기간이라면 정기결제 금액 지원 대신 SSG머니 3,900원 제공'', ''      스타벅스 자동 충전, 생활요금(통신요금, 아파트관리비 등) 정기결제 신청 및 이체 시 최대 1만원 청구 할인'', ''      2. 스타벅스 자동 충전 또는 생활요금 정기결제 신청(각 항목당 할인한도 5천원, 최대 1만원 할인)'', ''      정기결제 신청 후 카드 결제일에 따라 매출 발생 다음 달 또는 다다음 달 청구 할인 혜택 적용'', ''      단, 청구 할인 제공 일정은 당사 또는 신청인 사정에 의해 상이할 수 있음'', ''      3. 쏘카 1만원 할인쿠폰'', ''      쿠폰은 등록일 포함 30일간 이용 가능'', ''SSG.COM카드 Edition2는 실물카드 수령 전에도SSGPAY를 통해 바로 결제 가능합니다 (단, 실물카드 수령 전까지 SSG.COM 1회 결제 한도는 100만원으로 제한됨)'', ''SSG머니 최대 적립 12%에서 5%는 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능함'', ''카드 이용대금 연체 시 약정금리 + 연체가산금리 3%의 연체이자율이 적용됩니다. (회원별, 이용 상품별 차등적용 / 법정 최고금리 20% 이내) 단, 연체 발생시점에 약정금리가 없는 경우 아래와 같이 적용'', ''일시불 : 거래 발생시점 기준 최소 기간 (2개월)의 유이자 할부 약정금리 + 연체가산금리 3%'', ''무이자할부 : 거래발생시점 기준 동일한 할부 계약 기간의 유이자할부 약정금리 + 연체가산금리 3%'']'

# Row 6
$ws.Range("A6").Value = '더 강력해진 SSG.COM 삼성카드'
$ws.Range("B6").Value = 'SSG MONEY 최대 15% 적립 + 스마일클럽 월이용료 할인'
$ws.Range("C6").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000002385&recruitmentPath=SSG'
$ws.Range("D6").Value = 'SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2022-10-26'
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '2025-10-25'
$ws.Range("G6").Value = '[''이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''SSG머니 최대 15% 적립'', ''스마일클럽 월 이용료 3,900원 결제일 할인'', ''최대 15% SSG머니 적립 중 5%는 SSG.COM에서 제공하는 서비스로 자세한 내용은 SSG.COM 확인바람.'', ''01. SSG MONEY 최대 15% 적립!'', ''       이용실적 관계없이 적립 : 1 % + 이용실적 충족 시 적립 4% = 적립 가능한 최대 SSG머니 5%'', ''       카드혜택 + 스마일클럽 혜택'', ''       이용실적 관계없이 적립 1% + 이용실적 충족 시 적립 9% + 쓱배송/쓱배송 traders/새벽배송 상품 구매 시 5% = 적립 가능한 최대 SSG머니 15%'', ''최대 15% SSG MONEY 적립'', ''SSG.COM 삼성카드 최대 10% 적립+스마일클럽 5% 적립'', ''SSG.COM 삼성카드 최대 10% 적립(1% 적립+9% 추가 적립)'', ''1% 적립(전월 실적 조건 및 적립한도 없음)'', ''9% 추가 적립'', ''삼성카드 할인이 적용된 일시불 및 할부 이용금액은 제외됩니다.'', ''스마일클럽 5% 적립'', ''쓱배송/쓱배송 traders/새벽배송 상품 결제건에 한해 혜택을 받을 수 있습니다.'', ''02. 스마일클럽 월이용료 매월 3,900원 결제일 할인'', ''15,000원만 담아도 쓱 새벽배송 무료배송(SSG가입시)'', ''1쓱 새벽배송 트레이더 - SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''스마일클럽 단독 혜택 - 스타벅스 상품 전용 딜'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''카드 할인 혜택 자세히보기 (레이어팝업 열기)'', ''2022.10.26부터, SSG.COM 삼성카드 발급 시 스마일클럽에 자동 가입됩니다.'', ''G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우 혜택을 받을 수 없습니다.'', ''SSG.COM 삼성카드로 스마일클럽 월 이용료(3,900원) 정기결제 시 혜택이 제공됩니다.(월 1회)'', ''SSG.COM을 통한 스마일클럽 가입건에 한해 혜택이 제공됩니다.'', ''결제금액이 할인금액보다 적을 경우, 결제금액만큼 할인이 적용됩니다.(결제금액이 없는 경우 할인 대상에서 제외)'', ''01. 5만원 이상 결제 시 사용 가능한 4만 5천원 할인쿠폰 제공'', ''45,000원 할인쿠폰 - SSGPAY 바로결제 이용 시 사용 가능'', ''쿠폰 발급기간 : 2022.12.01 ~ 2022.12.31'', ''쿠폰 사용기간 : 2022.12.01 ~ 2022.12.31'', ''SSGPAY 바로결제 이용 이력이 없는 회원'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 50,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제하셔야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 통합 회원 본인 명의의 SSG.COM 삼성카드로 결제 시 사용 가능하며, 다른 부정적인 방법으로 사용한 경우에는 주문이 취소될 수 있습니다.'', ''쿠폰은 결제 화면에서 자동으로 적용됩니다.'', ''결제 화면에서 쿠폰 변경을 원할 경우 ‘쿠폰선택’을 눌러주세요.'', ''02. 삼성카드 결제일할인 행사일에는 최대 5% 즉시할인'', ''[카드혜택가 > SSG.COM 삼성카드 행사 더보기] 내 즉시할인가가 표시된 상품 결제시 적용(상품별 할인한도 30,000원)'', ''삼성카드 결제일할인 행사일에만 적용됩니다.'', ''삼성카드 할인이 적용된 결제건은 SSG머니 적립 대상에서 제외됩니다.'', ''03. 스마일클럽 고객이라면 9% 장바구니 쿠폰'', ''SSG.COM 삼성카드로 100,000원 이상 결제 시 사용 가능(아이디당 1회, 할인한도 20,000원)'', ''SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 100,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 결제 시 최적 할인에 의해 자동 설정되며, 직접 변경할 수 있 습니다.'', ''현금성 상품, 무형서비스 상품, 초특가 상품 등 일부 상품은 제외됩니다.(쿠폰 적용 제외 상품은 결제 화면에서 쿠폰이 보이지 않음)'', ''04. 스마일클럽 월 이용료 결제 카드로 SSG.COM 삼성카드 등록 시 SSG머니 3,900원 적립'', ''(1)행사기간 동안 SSG.COM 삼성카드를 통해 스마일클럽 자동 가입 시 SSG머니 3,900원 즉시 적립'', ''적립시점까지 스마일클럽 월 이용료 정기결제 수단에 SSG.COM 삼성카드를 등록해야 혜택을 받을 수 있습니다.'', ''본 상품 발급 및 SSGPAY앱 설치시 바로결제에 자동 등록 됩니다.'', ''SSG MONEY는 매월 1일~말일까지 매출전표가 접수된 금액에 대해 다음달 25일 SSG.COM 계정으로 적립됩니다.'', ''적립된 SSG MONEY는 SSGPAY 회원가입 후 조회 및 사용 가능합니다. 단, SSG.COM 에서는 SSG.COM만 가입해도 사용 가능합니다.'', ''연체이자율 : 회원별/이용상품별 정상이자율+3.0%p(최고 연 20.0%)'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!스마일클럽 가입하고 모든 혜택 누리세요'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우, 월 이용료 결제일할인 혜택을 받을 수 없습니다.'']'

# Row 7
$ws.Range("A7").Value = '연말 홈파티 메뉴 추천'
$ws.Range("B7").Value = '최대 30% 할인부터 페이백 이벤트까지'
$ws.Range("C7").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003406'
$ws.Range("D7").Value = 'Holiday at Home'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2022-12-19'
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '2022-12-31'
$ws.Range("G7").Value = '[''이벤트/쿠폰 > Holiday at Home'', ''스마일클럽'', ''공동 구매 이벤트'', '' SSG.COM 안의 작은 카페, Cafe at HOME 입니다. 12월 이벤트 토픽은 ‘홀리데이’에요. 홈파티 상품 할인부터 귀여운굿즈까지 영상으로 소개해드릴게요. '', ''EVENT01 홈파티 인기 상품 공구 이벤트'', ''지급된 SSG MONEY는 지급일로부터 30일간 사용 가능하며 기간 경과 후에는 소멸됩니다.'', ''본 이벤트는 당사의 사정에 따라 임의로 변경 혹은 조기 종료될 수 있습니다.'', ''새벽배송을 7만원 이상 주문하면, 카페앳홈의 귀여운 한정판 굿즈를 랜덤으로 알비백 안에 쏙 담아드립니다. 선착순 5백명이니, 서두르시면 좋아요!'', ''이벤트 기간 내 새벽배송 7만원 이상 주문 시 주문서에 사은품이 자동으로 포함됩니다.'', ''선착순 수량이 모두 소진되었을 시, 주문서에 사은품이 보이지 않습니다.'', ''      새벽배송 7만원 이상 주문 시, 선착순 증정'', ''주문서에 사은품이 보이지 않을 시, 선착순 수량이 모두 소진된 것임을 참고해주세요.'', ''본 이벤트는 당사 사정에 따라 임의로 변경 혹은 조기 종료될 수 있습니다.'']'

# Row 8
$ws.Range("A8").Value = 'W컨셉 베스트 브랜드 셀렉션'
$ws.Range("B8").Value = '13% 쿠폰+모두가 사랑한 브랜드를 한 자리에!'
$ws.Range("C8").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003454'
$ws.Range("D8").Value = '[12/26~1/1] 2022 W컨셉 BEST BRAND SELECTION'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2022-12-26'
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '2023-01-01'
$ws.Range("G8").Value = '[''이벤트/쿠폰 > [12/26~1/1] 2022 W컨셉 BEST BRAND SELECTION'', ''스마일클럽'', ''01) SPECIAL OFFER 13% 쿠폰 바로보기'', ''01) SPECIAL OFFER 13% 쿠폰 상품 쿠폰'', ''1만원 이상 구매시 최대 2만원 할인 -기간 내 ID당 4장 발급 가능-'', ''         쿠폰 발급 및 사용 기간'', ''         - 1만원 이상 상품 구매 시 13% 할인, 최대 2만원 (기간내 ID당 발급수량 4장)'', ''         쿠폰 적용 방법'', ''- 쿠폰은 적용상품 수량 1개에 대해 적용됩니다. (동일상품 복수 개 구매 시 1개 수량에만 적용)'']'

# Row 9
$ws.Range("A9").Value = '스마일클럽 12월의 가입 혜택'
$ws.Range("B9").Value = 'WELCOME 5천원 할인쿠폰 제공'
$ws.Range("C9").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003057'
$ws.Range("D9").Value = '[스마일클럽] 12월 내내 매일매일 스마일'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2022-12-05'
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = '2022-12-31'
$ws.Range("G9").Value = '[''이벤트/쿠폰 > [스마일클럽] 12월 내내 매일매일 스마일'', ''스마일클럽'', ''[스마일클럽] 12월 내내 매일매일 스마일'', '' 지금 스마일클럽 가입하고 모든 혜택 받기 '', ''한달 내내 막강한 스마일클럽 혜택'', ''스마일클럽 전용 혜택을 잡으세요!'', ''스마일클럽이 처음이라면 첫달 무료 + 5천원 할인 쿠폰 발급 (5,100원 이상 구매 시 사용 가능)'', ''본 쿠폰은 이벤트 기간 내 SSG.COM에서 스마일클럽에 신규 가입한 고객님에 한해 ID당 1회 발급됩니다.'', ''쿠폰은 가입 차주 금요일 이내 자동 지급되며, 앱푸시나 문자메시지 등을 통해 별도 안내 예정입니다. (단, 마케팅 정보 수신을 비동의 하신 경우 안내 발송이 제한됩니다.)'', ''본 쿠폰은 5,100원 이상 구매 시 5,000원 할인됩니다. (할인액 및 배송비 제외한 구매 금액 기준으로 쿠폰 적용됨)'', ''일부 상품 및 브랜드는 쿠폰 적용 제외될 수 있습니다.'', ''SSG.COM의 쿠폰은 결제 시 최적 할인에 의해 자동 설정되며, 직접 변경이 가능합니다.'', ''본 이벤트는 당사 사정으로 내용이 변경되거나 종료될 수 있습니다'', ''2023년 검은 토끼의 해 맞이 건강 프로젝트 스마일 클럽 전용 15% 쿠폰으로 건강을 선물해요!'', ''SSG.COM 카드 Edition2 매월 스마일클럽 가입비 100% 지원'', ''계속해서 이어질 스마일클럽 전용 혜택을 기대해주세요'']'

# Row 10
$ws.Range("A10").Value = '2023 스마일클럽 건강 프로젝트'
$ws.Range("B10").Value = '인기브랜드 15% 쿠폰 혜택'
$ws.Range("C10").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003386'
$ws.Range("D10").Value = '2023 스마일클럽 건강 프로젝트'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2022-12-19'
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = '2022-12-31'
$ws.Range("G10").Value = '[''이벤트/쿠폰 > 2023 스마일클럽 건강 프로젝트'', ''스마일클럽'', ''2023 스마일클럽 건강 프로젝트'', ''건강 TOP 브랜드 15% 쿠폰'', ''스마일클럽 회원 전용'', ''SmlieClub전용 15%상품쿠폰 : 5만원 이상 구매시 최대 1만 5천원 할인'', ''      쿠폰 발급 대상 및 발급 매수'', ''      쿠폰발급 및 사용 기간'', ''SmlieClub전용 12%장바구니쿠폰 : 3만원 이상 구매시 최대 1만원 할인'', ''쿠폰 대상 상품은 상품상세에서 다음 엠블럼을 확인해주세요'', ''스마일클럽 건강 프로젝트 최대 15%쿠폰 #연말연시선물 #새해새다짐 #탑브랜드'', ''쿠폰은 기간 내 ID당 1회 발급 가능합니다.'', ''본 쿠폰은 신세계몰, 신세계백화점, 이마트몰, 새벽배송, 트레이더스몰에 적용되는 쿠폰입니다.'', ''본 쿠폰은 장바구니 쿠폰으로 구매 건 내 해당 브랜드에 적용되며, 3만원 이상 주문 시 쿠폰이 적용됩니다.'', ''할인액 및 배송비 제외한 금액 기준으로 쿠폰이 적용됩니다. (에누리, 상품 할인 제외)'', "쿠폰 발급 후 ''My SSG > 쿠폰 > 보유쿠폰'' 에서 확인 가능합니다.", ''쿠폰발급 및 사용 기간 :  2022년 12월 19일(월) ~ 12월 31일(토)  '', ''겨울맞이 10% 즉시 할인중♥'', ''★15%즉시할인★ 종근당건강 BEST 행사 상품'', ''해피라운지★올인원 이뮨 액상 비타민 구매시 콜라겐 구미 증정★'', ''이마트 선물세트 할인 및 적립 상품권 증정까지 보러가기'']'

# Row 11
$ws.Range("A11").Value = '[월간캠핑쇼]  이타카 12/26(월) 7PM'
$ws.Range("B11").Value = '그랑데 스테이션 38만원+코베아 3웨이 구이바다 6.9만원대 핫딜'
$ws.Range("C11").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003558&siteNo=6005'
$ws.Range("D11").Value = '[쓱라이브] 이타카 12/26(월) 7PM'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2022-12-21'
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = '2022-12-31'
$ws.Range("G11").Value = '[''이벤트/쿠폰 > [쓱라이브] 이타카 12/26(월) 7PM'', ''스마일클럽'', '' - 사은품 지급 및 이벤트 혜택 제공'', '' - 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'

# Row 12
$ws.Range("A12").Value = '포레오 12/27(화) 8PM'
$ws.Range("B12").Value = '라이브중에만 추가 10% 쿠폰! 쏟아지는 구매인증 선물까지'
$ws.Range("C12").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003559&siteNo=6005'
$ws.Range("D12").Value = '[신세계라이브쇼핑] 포레오 12/27(화) 8PM'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2022-12-21'
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = '2022-12-31'
$ws.Range("G12").Value = '[''이벤트/쿠폰 > [신세계라이브쇼핑] 포레오 12/27(화) 8PM'', ''스마일클럽'', '' - 사은품 지급 및 이벤트 혜택 제공'', '' - 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'

# Row 13
$ws.Range("A13").Value = '12월 BEST 유아동 크리스마스 선물'
$ws.Range("B13").Value = '산타가 골라준 홀리데이 장난감'
$ws.Range("C13").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003103'
$ws.Range("D13").Value = '[12/01-31] 12월의 유아동 매거진 Little SSG X 스타일러'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2022-12-01'
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = '2022-12-31'
$ws.Range("G13").Value = '[''이벤트/쿠폰 > [12/01-31] 12월의 유아동 매거진 Little SSG X 스타일러'', ''스마일클럽'', ''"설레는 마음으로 크리스마스를 준비하는 시즌 스타일러맘과 함께하는 퀵 월간 쇼핑, ''리틀 쓱''과 함께해요!" 이런저런 고민 없이 클릭 한 번이면 육아템을 빠르게 구매할 수 있는 시대. 기존 ''유아동 전문관''에서 ''리틀 쓱(Little SSG)''으로 새롭게 선보이는 SSG닷컴(www.ssg.com)을 만나보세요. 계절에 꼭 필요한 각종 패션 아이템부터 아이와 집 안에서 보내는 시간을 풍요롭게 도와줄 교육 놀잇감, 친환경 육아용품 등 매달 <스타일러> 지면과 함께하는 유아동 기획전을 리틀 쓱 페이지에서도 소개합니다! ‘쓱 배송’은 물론 ‘타임 딜’ 등 수시로 열리는 다양한 이벤트 혜택까지, 현명한 스타일러맘이라면 지금 당장 리틀 쓱 페이지에 접속!'']'

# Row 14
$ws.Range("A14").Value = '12월 맘키즈 PLUS'
$ws.Range("B14").Value = '지금 할인 중! ~50% 할인 혜택'
$ws.Range("C14").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000001665'
$ws.Range("D14").Value = '이달의 맘키즈 PLUS'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2022-09-01'
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = '2999-12-13'
$ws.Range("G14").Value = '[''이벤트/쿠폰 > 이달의 맘키즈 PLUS'', ''스마일클럽'', ''맘키즈 ~40% 쿠폰상품'', ''지금 할인 중!'', ''※ 할인 금액은 상품 별로 상이하니 각 상품페이지를 꼭 참조하세요'', ''맘키즈 에누리 쿠폰'', ''행사대상 상품별 최대 40% 할인'', ''맘키즈 상품 에누리 쿠폰'', ''맘키즈 클럽 회원이라면 로그인 후 각 상품 상세페이지에서도 쿠폰을 받으실 수 있습니다.'', ''잠깐! 쿠폰이 적용되지 않는다면?'', ''맘키즈 에누리 쿠폰으로 상품별 최대 40%'', ''앱솔루트 유기농 궁 2 800g (NEO2 쓱배송, 그외지역 택배)'', ''콩순이 NEW콩순이냉장고 (N2쓱배송, 전국택배)'', ''자연나라 세종대왕김 국산 100% 도시락김 5g 36봉'', ''[맘키즈특가40%할인]서경한우 미경산 암소한우 1++등급  불고기 300g(냉장)'', ''헬로카봇 하이퍼캅스(N2쓱배송, 전국택배가능)'']'

Write-Host "Sheet updated"